$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, reusing the same header style as the
# existing header cells (e.g. G1) by copying formats over.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the corresponding data value in H2 (row 2), matching the diff.
$ws.Range("H2").Value = 0
